# Regenerate save_data column G ("K") values for rows 2-24.
# These values were recomputed (std/mean, s_vals) and rewritten to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 7
    4  = 5
    5  = 4
    6  = 11
    7  = 7
    8  = 9
    9  = 5
    10 = 9
    11 = 3
    12 = 7
    13 = 4
    14 = 4
    15 = 5
    16 = 9
    17 = 4
    18 = 4
    19 = 5
    20 = 5
    21 = 4
    22 = 2
    23 = 3
    24 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
